$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.376.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "'1.827.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'313.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.4480"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.81%  "
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("D9").Value = "'0.07501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.12%  "
$ws.Range("D10").Value = "'0.8972"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.98%  "
$ws.Range("D11").Value = "'21.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").Value = "'1.811.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'6.783"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "'94.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.03%  "
$ws.Range("D15").Value = "'5.413"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "'0.07122"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'0.000008822"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "'0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'15.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").Value = "'27.395.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Value = "'5.293"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.49%  "
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "'2.055.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'2.003"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").Value = "'2.500"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.90%  "
$ws.Range("D27").Value = "'151.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "'18.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").Value = "'5.388"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("D30").Value = "'118.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'0.08856"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "'0.7806"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.58%  "
$ws.Range("D33").Value = "'1.204"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").Value = "'4.592"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("D35").Value = "'2.887"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "'0.9989"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "'1.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("E38").Value = "  +2.70%  "
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").Value = "'7.408"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("D41").Value = "'0.5354"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.56%  "
$ws.Range("D42").Value = "'0.1733"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").Value = "'2.857"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'2.294"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +18.04%  "
$ws.Range("D45").Value = "'8.829"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("D46").Value = "'0.5161"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.48%  "
$ws.Range("D47").Value = "'10.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("D48").Value = "'106.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("D50").Value = "'0.9991"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'0.06388"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
